# Planning previsionnel.xlsx - update
#  - advance the "selected period" spinner from 2 to 3 (N3, linked to
#    the "période_sélectionnée" defined name)
#  - fill in the actual start/duration (E/F) and percentage achieved (G)
#    for the remaining tasks in the Gantt table
#  - change the view (zoom + selection) to what was on screen when saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Spinner-linked cell: période_sélectionnée (projet!$N$3) goes from 2 to 3.
$ws.Range("N3").Value = 3

# Gantt rows 11-12 already have "réel" dates; only their completion % moves to 100%.
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 1

# Rows 13-27: fill in actual start (E) / duration (F), and mark as complete (G = 1)
# except row 27, whose completion stays at 0.
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 1

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1

$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 1

$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 4
$ws.Range("G19").Value = 1

$ws.Range("E20").Value = 4
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 1

$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 3
$ws.Range("G21").Value = 1

$ws.Range("E22").Value = 6
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1

$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1

$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 1

$ws.Range("E25").Value = 7
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 1

$ws.Range("E26").Value = 8
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1

$ws.Range("E27").Value = 9
$ws.Range("F27").Value = 1
# G27 is left at 0 - unchanged

# View state: zoomed out to 80% and scrolled/selected A26.
$excel.ActiveWindow.Zoom = 80
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A26").Select()
